$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.672.01'
$ws.Range("E2").Value = '  +0.16%  '

$ws.Range("D3").Value = '1.843.46'
$ws.Range("E3").Value = '  -0.24%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.30'
$ws.Range("E5").Value = '  +0.89%  '

$ws.Range("E6").Value = '  +0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4316'
$ws.Range("E7").Value = '  +0.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3707'
$ws.Range("E8").Value = '  +2.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07335'
$ws.Range("E9").Value = '  +0.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8780'
$ws.Range("E10").Value = '  +0.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.04'
$ws.Range("E11").Value = '  +1.69%  '

$ws.Range("D12").Value = '1.936.32'
$ws.Range("E12").Value = '  +4.79%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.474'
$ws.Range("E13").Value = '  +2.75%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.599'
$ws.Range("E14").Value = '  +1.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06955'
$ws.Range("E15").Value = '  +0.63%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.004'
$ws.Range("E16").Value = '  +0.34%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.12'
$ws.Range("E17").Value = '  +1.42%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009054'
$ws.Range("E18").Value = '  +0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.002'
$ws.Range("E19").Value = '  +0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.56'
$ws.Range("E20").Value = '  +1.53%  '

$ws.Range("D21").Value = '27.953.89'
$ws.Range("E21").Value = '  +1.11%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.121'
$ws.Range("E22").Value = '  +3.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.99'
$ws.Range("E23").Value = '  +5.83%  '

$ws.Range("D24").Value = '2.143.76'
$ws.Range("E24").Value = '  +3.29%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.990'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.11'
$ws.Range("E26").Value = '  -0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.93'
$ws.Range("E27").Value = '  +0.76%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.319'
$ws.Range("E28").Value = '  +0.58%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.85'
$ws.Range("E29").Value = '  -4.69%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.876'
$ws.Range("E30").Value = '  +1.53%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08915'
$ws.Range("E31").Value = '  +0.14%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7890'
$ws.Range("E32").Value = '  +3.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.619'
$ws.Range("E33").Value = '  +1.50%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.174'
$ws.Range("E34").Value = '  +6.43%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.961'
$ws.Range("E35").Value = '  -0.54%  '

$ws.Range("E36").Value = '  +0.14%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05443'
$ws.Range("E37").Value = '  +0.58%  '

$ws.Range("E38").Value = '  +1.58%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01962'
$ws.Range("E39").Value = '  +1.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.843'
$ws.Range("E40").Value = '  +0.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5176'
$ws.Range("E41").Value = '  +1.93%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1692'
$ws.Range("E42").Value = '  +2.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.809'
$ws.Range("E43").Value = '  +0.74%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.649'
$ws.Range("E44").Value = '  +3.31%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.66'
$ws.Range("E45").Value = '  +2.77%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4781'
$ws.Range("E46").Value = '  +2.19%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '106.62'
$ws.Range("E47").Value = '  +1.37%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06549'
$ws.Range("E48").Value = '  +0.02%  '

$ws.Range("E49").Value = '  +0.10%  '

$ws.Range("E50").Value = '  +2.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.846'
$ws.Range("E51").Value = '  +5.93%  '
